$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number-format (date) of the previous row's date cell so the new
# cell reuses the existing style instead of creating a duplicate one.
$ws.Range("A89").Copy()
$ws.Range("A90").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row of data for 2022-06-30 (date serial 44742)
$ws.Range("A90").Value = 44742
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 330327
$ws.Range("D90").Value = 6358
$ws.Range("E90").Value = 253
$ws.Range("F90").Value = 0

$ws.Range("F90").Select()
